# Apply the canonical terminology URL / date update described by the diff.
$wb = $excel.ActiveWorkbook

# Sheet "Metadata": update the Date value (row 8, column B)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# Sheet "Include #0": update the System URI for TRE-R39-Competence (row 4, column B)
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R39-Competence/FHIR/TRE-R39-Competence"

# Sheet "Include #1": update the System URI for TRE-R01-EnsembleSavoirFaire-CISIS (row 4, column B)
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R01-EnsembleSavoirFaire-CISIS/FHIR/TRE-R01-EnsembleSavoirFaire-CISIS"
